# NSMB - 8-F1 up to boss fight.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("V4")

# --- Row 111: drop the "Frame ruled out due to moving logs" note from G111 ---
$ws.Cells.Item(111, 7).ClearContents()

# --- Row 112: G112 note changes from "-" to "Frame ruled out due to moving logs" ---
$ws.Cells.Item(112, 7).Value = "Frame ruled out due to moving logs"

# --- Row 113: G113 stays "-" (unchanged) ---

# --- Row 114: B114 corrected 34776 -> 34767, G114 (545) removed ---
$ws.Cells.Item(114, 2).Value = 34767
$ws.Cells.Item(114, 7).ClearContents()

# --- Insert four new rows (new checkpoints) before the old row 115 ---
$ws.Rows("115:118").Insert()

# Copy formatting (style) from row 114 A:D down into the newly inserted rows
$ws.Range("A114:D114").Copy()
$ws.Range("A115:D118").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 115
$ws.Cells.Item(115, 1).Value = "Wall bounce dust"
$ws.Cells.Item(115, 2).Value = 34904
$ws.Cells.Item(115, 3).Value = 40799
$ws.Cells.Item(115, 4).Formula = "=IF(B115 >  0,C115-B115, 0)"

# Row 116
$ws.Cells.Item(116, 1).Value = "Wall bounce dust"
$ws.Cells.Item(116, 2).Value = 34950
$ws.Cells.Item(116, 3).Value = 40846
$ws.Cells.Item(116, 4).Formula = "=IF(B116 >  0,C116-B116, 0)"

# Row 117
$ws.Cells.Item(117, 1).Value = "Wall bounce dust"
$ws.Cells.Item(117, 2).Value = 35025
$ws.Cells.Item(117, 3).Value = 40921
$ws.Cells.Item(117, 4).Formula = "=IF(B117 >  0,C117-B117, 0)"
$ws.Cells.Item(117, 8).Value = 41175
$ws.Cells.Item(117, 10).Formula = "=35282-74"

# Row 118
$ws.Cells.Item(118, 1).Value = "Checkpoint 626"
$ws.Cells.Item(118, 2).Value = 35282
$ws.Cells.Item(118, 3).Value = 41175
$ws.Cells.Item(118, 4).Formula = "=IF(B118 >  0,C118-B118, 0)"
$ws.Cells.Item(118, 8).Value = 41101

# Row 119 (was the old row 115 - "Enter door") now also gets a B value and H formula
$ws.Cells.Item(119, 2).Value = 35587
$ws.Cells.Item(119, 4).Formula = "=IF(B119 >  0,C119-B119, 0)"
$ws.Cells.Item(119, 8).Formula = "=H117-H118"

# --- View: keep header row frozen, move viewport down, select B120 ---
$win = $excel.ActiveWindow
$win.FreezePanes = $false
$null = $ws.Range("A2").Select()
$win.FreezePanes = $true
$null = $ws.Range("B120").Select()
